# Generate Report for Archive
# The localization status for the two handed-off files moved on from
# "Ready for handoff" to "In Translation" - update every sheet that
# surfaces the Status column (the Overview roll-up as well as each
# per-locale detail sheet), then let Excel re-fit the Status columns
# now that the text is shorter.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: Status is mirrored per-locale in columns E (zh-cn) and F (de-de)
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

# --- Per-locale detail sheets: Status lives in column C
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Re-fit the Status columns now that the text is narrower
$overview.Columns.Item(5).AutoFit()
$overview.Columns.Item(6).AutoFit()
$zhcn.Columns.Item(3).AutoFit()
$dede.Columns.Item(3).AutoFit()
